# Apply "stock updated by raj time 12:59" changes to the Stock Summary workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PATRIKA 24-25")

# Update the date-range label (appears twice: A3 and B6)
$ws.Range("A3").Value = "1-Jul-2024 to 28-Dec-2024"
$ws.Range("B6").Value = "1-Jul-2024 to 28-Dec-2024"

# Row -> (Quantity (B), Value (D))
$updates = @{
    24  = @(57,      119.7)
    25  = @(111.5,   234.15)
    45  = @(242.5,   266.75)
    46  = @(370,     407)
    50  = @(250,     275)
    67  = @(105.5,   94.95)
    68  = @(359,     556.45)
    77  = @(484.5,   557.18)
    78  = @(463.5,   533.03)
    161 = @(88,      334.4)
    182 = @(80.5,    305.9)
    185 = @(15,      61.68)
    207 = @(71,      404.7)
    259 = @(20.5,    147.14)
    364 = @(66,      178.2)
    383 = @(37.06,   129.71)
    412 = @(43,      306.59)
    433 = @(20,      60)
    441 = @(80,      136.8)
    480 = @(34.5,    172.5)
    527 = @(48.5,    196.43)
    531 = @(57.5,    80.5)
    538 = @(550,     440)
    540 = @(400,     320)
    553 = @(40140.72, 111589.92)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}
